# Regenerate save_data to use K instead of Strike# (column G values),
# recalculated std/mean based s_vals written back as literal values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 5
    4  = 3
    5  = 7
    6  = 3
    7  = 2
    8  = 4
    9  = 2
    10 = 3
    11 = 9
    12 = 5
    13 = 4
    14 = 4
    15 = 3
    16 = 3
    17 = 5
    18 = 6
    19 = 5
    20 = 2
    21 = 6
    22 = 4
    23 = 6
    24 = 5
    25 = 8
    26 = 2
    27 = 4
    28 = 1
    29 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
